# Juno: check in to OLPRODLOC.
# Rename the worksheet and translate the regional column headers to
# Simplified Chinese.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (tab) from "Sales report" to the Chinese equivalent.
$ws.Name = "销售报表"

# Translate region header labels in row 1 (Year-Quarter / Mountain stay
# in English; the other four regions are localized).
$ws.Range("B1").Value = "中西部"
$ws.Range("D1").Value = "东北"
$ws.Range("E1").Value = "南部"
$ws.Range("F1").Value = "东南部"
$ws.Range("G1").Value = "西部"
